$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper block builder. Each "table" block consists of:
#   row N   : merged title row (A:D), bold 14pt centered
#   row N+1 : column-header row (bold, light-grey fill)
#   row N+2 : value/label row (PK/FK markers, grey-ish "Output" look)
# We replicate formatting by copying the already-styled Admin Table block
# (rows 3-5) which Excel itself used as the template when the author
# duplicated it for the new tables.
# ---------------------------------------------------------------------------

function Copy-BlockFormat($titleRow, $numCols) {
    # Copy the 4-column (A:D) formats from the existing Admin Table block.
    $ws.Range("A3:D5").Copy()
    $ws.Range("A" + $titleRow).PasteSpecial(-4122)

    if ($numCols -gt 4) {
        # Extend header/label formatting rightwards using column D as the template
        # (same header style + same label/value style) for the extra columns.
        $headerRow = $titleRow + 1
        $labelRow = $titleRow + 2
        $srcHeader = "D4"
        $srcLabel = "D5"
        for ($c = 5; $c -le $numCols; $c++) {
            $colLetter = [char](64 + $c)
            $ws.Range("D4:D5").Copy()
            $ws.Range($colLetter + $headerRow).PasteSpecial(-4122)
        }
    }
}

# --- Classroom Table (rows 18-20), columns A-E -----------------------------
Copy-BlockFormat 18 5
$ws.Range("A18").Value = "Classroom Table"
$ws.Range("A18:D18").Merge()

$ws.Range("A19").Value = "id"
$ws.Range("B19").Value = "name"
$ws.Range("C19").Value = "created_by"
$ws.Range("D19").Value = "class_code"
$ws.Range("E19").Value = "created_date"

$ws.Range("A20").Value = "PK"
$ws.Range("D20").Value = "unique"

# --- Quiz Table (rows 23-25), columns A-G -----------------------------------
Copy-BlockFormat 23 7
$ws.Range("A23").Value = "Quiz Table"
$ws.Range("A23:D23").Merge()

$ws.Range("A24").Value = "id"
$ws.Range("B24").Value = "classroom_id"
$ws.Range("C24").Value = "quiz_type"
$ws.Range("D24").Value = "total_marks"
$ws.Range("E24").Value = "total_time"
$ws.Range("F24").Value = "quiz_timeline"
$ws.Range("G24").Value = "quiz_date"

$ws.Range("A25").Value = "PK"
$ws.Range("B25").Value = "FK"

# --- Student_Quiz_record Table (rows 28-30), columns A-G --------------------
Copy-BlockFormat 28 7
$ws.Range("A28").Value = "Student_Quiz_record Table"
$ws.Range("A28:D28").Merge()

$ws.Range("A29").Value = "quiz_id"
$ws.Range("B29").Value = "classroom_id"
$ws.Range("C29").Value = "student_id"
$ws.Range("D29").Value = "correct_answers"
$ws.Range("E29").Value = "wrong_answers"
$ws.Range("F29").Value = "obtained_marks"
$ws.Range("G29").Value = "obtained_grade"

$ws.Range("A30").Value = "FK"
$ws.Range("B30").Value = "FK"
$ws.Range("C30").Value = "FK"

# --- Quiz_question_record Table (rows 33-35), columns A-G -------------------
Copy-BlockFormat 33 7
$ws.Range("A33").Value = "Quiz_question_record Table"
$ws.Range("A33:D33").Merge()

$ws.Range("A34").Value = "quiz_id"
$ws.Range("B34").Value = "classroom_id"
$ws.Range("C34").Value = "question_type"
$ws.Range("D34").Value = "question_mark"
$ws.Range("E34").Value = "correct_answer"
$ws.Range("F34").Value = "question_definition"
$ws.Range("G34").Value = "question_number"

$ws.Range("A35").Value = "FK"
$ws.Range("B35").Value = "FK"

# --- Column widths (B:G), matching the author's manual resize ---------------
$ws.Columns.Item(2).ColumnWidth = 13.5546875 - (5.0/6.0)
$ws.Columns.Item(3).ColumnWidth = 17.109375 - (5.0/6.0)
$ws.Columns.Item(4).ColumnWidth = 16.6640625 - (5.0/6.0)
$ws.Columns.Item(5).ColumnWidth = 16.5546875 - (5.0/6.0)
$ws.Columns.Item(6).ColumnWidth = 19.44140625 - (5.0/6.0)
$ws.Columns.Item(7).ColumnWidth = 17.21875 - (5.0/6.0)

# --- View state: active selection / scroll position -------------------------
try { $excel.ActiveWindow.ScrollRow = 13 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("F43").Select()

Write-Host "edit complete"
